# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.237.65"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "1.645.81"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.84"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0636"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.93"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.873.48"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.29"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "1.656.23"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.547"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.23"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "26.231.51"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.12"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.04"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.31"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.80"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.96"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0501"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.910"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "1.133.58"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.553"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0157"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.53"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.27"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.800"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.782.50"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.39"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("E47").Value = "  +3.90%  "
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.418"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.67"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0967"
$ws.Range("E51").Value = "  +1.65%  "
